$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "PLC2 Global" sheet - add two new rows (16 & 17) at the bottom
#    errorSystem / D3001 / WORD / 0 / "Transfer PLC system error to PLC1"
#    errorIO     / D3002 / WORD / 0 / "Transfer PLC I/O error to PLC1"
# ------------------------------------------------------------------
$wsPlc2 = $wb.Worksheets.Item("PLC2 Global")

$wsPlc2.Range("A16").Value = "errorSystem"
$wsPlc2.Range("B16").Value = "D3001"
$wsPlc2.Range("A17").Value = "errorIO"
$wsPlc2.Range("B17").Value = "D3002"

$wsPlc2.Range("E16").Value = "Transfer PLC system error to PLC1"
$wsPlc2.Range("E17").Value = "Transfer PLC I/O error to PLC1"

# Column C (type) on existing rows uses the bold/left style (s=8) -
# copy that formatting from the row above so the new cells reuse the
# same style record instead of minting a new one.
$wsPlc2.Range("C15").Copy()
$wsPlc2.Range("C16").PasteSpecial(-4122)
$wsPlc2.Range("C17").PasteSpecial(-4122)
$wsPlc2.Range("C16").Value = "WORD"
$wsPlc2.Range("C17").Value = "WORD"

$wsPlc2.Range("D16").Value = 0
$wsPlc2.Range("D17").Value = 0

# ------------------------------------------------------------------
# 2) "IO Mapping" sheet - tag existing row 163 with a new comment and
#    append two new rows (164 & 165) for the PLC2 transfer errors
# ------------------------------------------------------------------
$wsIo = $wb.Worksheets.Item("IO Mapping")

$wsIo.Range("A164").Value = "errorPLC2system"
$wsIo.Range("A165").Value = "errorPLC2io"

$wsIo.Range("B164").Value = "D27001"
$wsIo.Range("B165").Value = "D27002"

$wsIo.Range("F163").Value = "PLC2 VFD error"
$wsIo.Range("F164").Value = "PLC2 system error"
$wsIo.Range("F165").Value = "PLC2 IO error"

$wsIo.Range("C164").Value = "WORD"
$wsIo.Range("C165").Value = "WORD"
$wsIo.Range("D164").Value = 0
$wsIo.Range("D165").Value = 0

# ------------------------------------------------------------------
# 3) "Constants" sheet - insert a new row 69 for the PLC1/PLC2
#    transfer status flag, shifting the rest of the table down by one
# ------------------------------------------------------------------
$wsConst = $wb.Worksheets.Item("Constants")

$wsConst.Rows.Item(69).Insert()

$wsConst.Range("B69").Value = "stat_plc_transfer"
$wsConst.Range("G69").Value = "Transfer between PLC1 and PLC2"

$wsConst.Range("C69").Value = 1
$wsConst.Range("D69").Value = "BOOL"
$wsConst.Range("E69").Value = $true
$wsConst.Range("F69").Value = "x"

# ------------------------------------------------------------------
# 4) View / selection bookkeeping - mirror where the author left the
#    cursor on each sheet (best effort; Constants ends up the active
#    tab again, matching the saved file).
# ------------------------------------------------------------------
$wsIo.Select()
$excel.ActiveWindow.ScrollRow = 151
$excel.ActiveWindow.ScrollColumn = 1
$wsIo.Range("H163").Select()

$wsPlc2.Select()
$wsPlc2.Range("A16").Select()

$wsConst.Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 2
$wsConst.Range("G70").Select()
